$d = $word.ActiveDocument

# Every bold heading/label run in the document (section headings, the
# "Contact Information:" line, the summary blurb, etc.) needs an explicit
# Italic=False applied, EXCEPT the very first paragraph (the "Ethan
# Ransberger" title). Skip paragraph 1 and, for every other paragraph whose
# whole range is already bold, force Italic off across the paragraph's text
# (excluding the trailing paragraph mark, so we don't touch w:pPr/w:rPr).
$count = $d.Paragraphs.Count
for ($i = 2; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $pRange = $p.Range
    if ($pRange.Font.Bold -eq -1) {
        $endPos = $pRange.End - 1
        if ($endPos -gt $pRange.Start) {
            $r = $d.Range($pRange.Start, $endPos)
            $r.Font.Italic = 0
        }
    }
}
